# Update gh-pages to output generated at 456a3b4
# Apply the same set of value updates to both the "展览" sheet and the
# "全部类型" sheet (which duplicates the same rows), matching the diff.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    if ($name -eq "展览") {
        $ws.Range("F2").Value = 464
        $ws.Range("F3").Value = 5564
        $ws.Range("F6").Value = 83
        $ws.Range("F8").Value = 51
        $ws.Range("F9").Value = 526
        $ws.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202409/mDTW4lHL1727235917704.jpeg"
        $ws.Range("F10").Value = 19
    }
    elseif ($name -eq "全部类型") {
        $ws.Range("F2").Value = 464
        $ws.Range("F3").Value = 5564
        $ws.Range("F7").Value = 83
        $ws.Range("F10").Value = 51
        $ws.Range("F11").Value = 526
        $ws.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202409/mDTW4lHL1727235917704.jpeg"
        $ws.Range("F12").Value = 19
    }
}
